# Fix survey and question/answer upload
#
# The "Survey" sheet's header row (C1:G1) held the raw answer-option
# numbers (1..5) instead of the option labels used elsewhere in the
# workbook ("A1".."A5"). Re-enter them as text labels so the upload
# matches the Options sheet's answer codes.
#
# The "count" column (H) on rows 2 and 4 was being read/uploaded as a
# number, which strips any leading zeros / causes mis-typing downstream;
# format those two cells as Text.
#
# Finally, re-select the Survey tab (it's the primary sheet that should
# be shown/active when the file is opened) instead of Options.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Survey")

$ws.Range("C1").Value = "A1"
$ws.Range("D1").Value = "A2"
$ws.Range("E1").Value = "A3"
$ws.Range("F1").Value = "A4"
$ws.Range("G1").Value = "A5"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H4").NumberFormat = "@"

$ws.Activate()
